# Restructure the LOM3106 summary sheet to match the updated course data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet shrinks from 25 to 23 rows; drop the two trailing rows first
# (bottom-up, so row numbers above stay stable while we still reference them).
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(24).Delete()

# Row 10
$ws.Range("B10").Value = "3480026 - João Paulo Pascon"
$ws.Range("C10").Value = "3480026 - João Paulo Pascon"

# Row 13
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "01/01/2022"
$ws.Range("C13").Value = "01/01/2022"
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Image processing in materialography; Adjusting empirical equations; Interatomic potentials and classical molecular dynamics; Description of nucleation and growth kinetics; Finite Element Method; Monte Carlo methods; Grain growth; Calculation of phase diagrams."
$ws.Range("C14").Value = "Image processing in materialography; Adjusting empirical equations; Interatomic potentials and classical molecular dynamics; Description of nucleation and growth kinetics; Finite Element Method; Monte Carlo methods; Grain growth; Calculation of phase diagrams."
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "3480026 - João Paulo Pascon"
$ws.Range("C15").Value = "3480026 - João Paulo Pascon"
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "- Image treatment: resolution, definition, contrast, saturation; use of automated techniques for determining particle size and distribution.- Proposition and fit of empirical equations to results of experimental measures: the various proposals for relationships for plastic deformation and hardening.- Interatomic potentials and the classical molecular dynamics method; simulation of solidification of a pure metal.- Nucleation and growth kinetics: the Johnson-Mehl-Avrami-Kolmogorov (JMAK) equation and its computational application.- Finite element method: study of the stress state of materials under mechanical loads; simulation of heat transfer applied to heat treatments.- Monte Carlo method applied to the ferro-paramagnetic transition and to grain growth kinetics- Calculation of phase diagrams: free energy curves, the CALPHAD method; Thermo-Calc and Dictra."
$ws.Range("C16").Value = "- Image treatment: resolution, definition, contrast, saturation; use of automated techniques for determining particle size and distribution.- Proposition and fit of empirical equations to results of experimental measures: the various proposals for relationships for plastic deformation and hardening.- Interatomic potentials and the classical molecular dynamics method; simulation of solidification of a pure metal.- Nucleation and growth kinetics: the Johnson-Mehl-Avrami-Kolmogorov (JMAK) equation and its computational application.- Finite element method: study of the stress state of materials under mechanical loads; simulation of heat transfer applied to heat treatments.- Monte Carlo method applied to the ferro-paramagnetic transition and to grain growth kinetics- Calculation of phase diagrams: free energy curves, the CALPHAD method; Thermo-Calc and Dictra."
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Rows.Item(17).EntireRow.AutoFit()

# Row 18
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados. Trabalho baseado em Projeto"
$ws.Range("C19").Value = "Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados. Trabalho baseado em Projeto"
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média aritmética de trabalhos propostos ao longo do curso (60%) e do Trabalho final em grupo (40%)."
$ws.Range("C20").Value = "Média aritmética de trabalhos propostos ao longo do curso (60%) e do Trabalho final em grupo (40%)."

# Row 21
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Não haverá exame de recuperação."
$ws.Range("C21").Value = "Não haverá exame de recuperação."
$ws.Rows.Item(21).RowHeight = 120

# Row 22
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Rows.Item(22).EntireRow.AutoFit()

# Row 23
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOM3013 -  Ciência dos Materiais  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOM3013 -  Ciência dos Materiais  (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30
